$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are written as text (matching the original inline-string cells),
# using a Text number format during the write so Excel does not coerce
# numeric-looking or percent-looking strings into real numbers, then the
# cell style is reset back to Normal to match the original formatting.
$cellUpdates = @(
    @{ Addr = 'D2'; Value = '311.01' }
    @{ Addr = 'E2'; Value = '-1.61%' }
    @{ Addr = 'D3'; Value = '38.20' }
    @{ Addr = 'E3'; Value = '-3.04%' }
    @{ Addr = 'D4'; Value = '5.067' }
    @{ Addr = 'E4'; Value = '-1.17%' }
    @{ Addr = 'D5'; Value = '0.07756' }
    @{ Addr = 'E5'; Value = '-5.11%' }
    @{ Addr = 'D6'; Value = '4.355' }
    @{ Addr = 'E6'; Value = '-0.33%' }
    @{ Addr = 'D7'; Value = '1.889' }
    @{ Addr = 'E7'; Value = '-4.92%' }
    @{ Addr = 'D8'; Value = '8.192' }
    @{ Addr = 'D9'; Value = '0.9221' }
    @{ Addr = 'E9'; Value = '-1.64%' }
    @{ Addr = 'D10'; Value = '0.1232' }
    @{ Addr = 'E10'; Value = '-5.21%' }
    @{ Addr = 'D11'; Value = '0.1878' }
    @{ Addr = 'E11'; Value = '-4.89%' }
    @{ Addr = 'D12'; Value = '0.08823' }
    @{ Addr = 'E12'; Value = '-2.73%' }
    @{ Addr = 'D13'; Value = '0.03440' }
    @{ Addr = 'E13'; Value = '-1.68%' }
    @{ Addr = 'D14'; Value = '0.09706' }
    @{ Addr = 'E14'; Value = '-0.26%' }
    @{ Addr = 'D15'; Value = '0.001370' }
    @{ Addr = 'E15'; Value = '-2.70%' }
    @{ Addr = 'D16'; Value = '0.006050' }
    @{ Addr = 'E16'; Value = '-0.11%' }
    @{ Addr = 'D17'; Value = '3.561' }
    @{ Addr = 'E17'; Value = '-1.98%' }
    @{ Addr = 'E18'; Value = '-6.48%' }
    @{ Addr = 'D19'; Value = '0.3409' }
    @{ Addr = 'E19'; Value = '-2.32%' }
    @{ Addr = 'D20'; Value = '5.030' }
    @{ Addr = 'E20'; Value = '1.36%' }
    @{ Addr = 'E21'; Value = '-2.60%' }
    @{ Addr = 'D22'; Value = '0.2619' }
    @{ Addr = 'E22'; Value = '1.52%' }
    @{ Addr = 'E23'; Value = '5,593.14%' }
    @{ Addr = 'D24'; Value = '0.04393' }
    @{ Addr = 'E24'; Value = '0.85%' }
    @{ Addr = 'D25'; Value = '0.001210' }
    @{ Addr = 'E25'; Value = '-2.48%' }
    @{ Addr = 'D26'; Value = '0.004254' }
    @{ Addr = 'E26'; Value = '-10.69%' }
    @{ Addr = 'E27'; Value = '-65.28%' }
    @{ Addr = 'D39'; Value = '0.02133' }
    @{ Addr = 'E39'; Value = '-3.59%' }
    @{ Addr = 'D40'; Value = '0.05021' }
    @{ Addr = 'E40'; Value = '-3.21%' }
    @{ Addr = 'D41'; Value = '0.007843' }
    @{ Addr = 'E41'; Value = '1.17%' }
    @{ Addr = 'D42'; Value = '0.009995' }
    @{ Addr = 'E42'; Value = '-3.50%' }
    @{ Addr = 'D43'; Value = '0.1342' }
    @{ Addr = 'E43'; Value = '-4.24%' }
    @{ Addr = 'D44'; Value = '0.001993' }
    @{ Addr = 'E44'; Value = '-5.14%' }
    @{ Addr = 'D45'; Value = '0.009704' }
    @{ Addr = 'E45'; Value = '4.55%' }
    @{ Addr = 'D46'; Value = '0.00006461' }
    @{ Addr = 'E46'; Value = '-6.96%' }
    @{ Addr = 'D47'; Value = '0.00000000750' }
    @{ Addr = 'E47'; Value = '-0.02%' }
    @{ Addr = 'D48'; Value = '0.003212' }
    @{ Addr = 'E48'; Value = '11.34%' }
    @{ Addr = 'E49'; Value = '-0.12%' }
    @{ Addr = 'D50'; Value = '0.00002101' }
    @{ Addr = 'E50'; Value = '-0.02%' }
    @{ Addr = 'D51'; Value = '0.0002001' }
    @{ Addr = 'E51'; Value = '-0.02%' }
)

foreach ($update in $cellUpdates) {
    $cell = $ws.Range($update.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = "Normal"
}
